# Update product listing rows with new data per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "5.99`" Смартфон TP-Link Neffos A5 16 ГБ зеленый"
$ws.Cells.Item(2, 2).Value = "'1645574"
$ws.Cells.Item(2, 4).Value = "https://s.technopoint.ru/thumb/st1/fit/800/650/9710105c98ee88b56ff0a7d79e1e8a05/6557ae69ab8ab3b7c6fc2918952a6d7105c54c7c0ea40dfc6b2dcaef45b1dbdb.png"

$ws.Cells.Item(3, 1).Value = "5.99`" Смартфон TP-Link Neffos A5 16 ГБ серый"
$ws.Cells.Item(3, 2).Value = "'1645572"
$ws.Cells.Item(3, 4).Value = "https://s.technopoint.ru/thumb/st1/fit/800/650/06e76fe70d7e293fb5e8627c2a4bdfe8/c958b42f59f71fce27c8bbd954be0361804ce47e23933231a127f45e98cc9a9e.png"

$ws.Cells.Item(4, 1).Value = "6.35`" Смартфон bright & quick BQ 6424L MAGIC O 32 ГБ красный"
$ws.Cells.Item(4, 2).Value = "'1646117"
$ws.Cells.Item(4, 3).Value = "'7450"
$ws.Cells.Item(4, 4).Value = "https://s.technopoint.ru/thumb/st4/fit/800/650/1223cae5f4db3cac95031a8da1ffd25d/92c5ed025c2f5179da87f639710916165f9c2fbe2e49490df61b06451e54388e.jpg"

$ws.Cells.Item(5, 1).Value = "5.7`" Смартфон Samsung Galaxy A01 16 ГБ черный"
$ws.Cells.Item(5, 2).Value = "'1623527"
$ws.Cells.Item(5, 3).Value = "'7999"
$ws.Cells.Item(5, 4).Value = "https://s.technopoint.ru/thumb/st1/fit/800/650/75a04bed348bd860b23d44282764dd9a/b9a56bbadd4cf887d946ca418544bac52432f7db4878b09aaf65bc23a91b0e4b.jpg"

$ws.Cells.Item(6, 1).Value = "6.52`" Смартфон realme C3 64 ГБ синий"
$ws.Cells.Item(6, 2).Value = "'1641173"
$ws.Cells.Item(6, 3).Value = "'9999"
$ws.Cells.Item(6, 4).Value = "https://s.technopoint.ru/thumb/st1/fit/800/650/e4ada0fd3b6d661d68c6b977f00e4bb0/2ee2e247fef5a4be9556895a05d2af5b2bdb3b6bdd66af02192b128e18196a7b.jpg"

$ws.Cells.Item(7, 1).Value = "6.5`" Смартфон OPPO A31 64 ГБ белый"
$ws.Cells.Item(7, 2).Value = "'1642439"
$ws.Cells.Item(7, 3).Value = "'11999"
$ws.Cells.Item(7, 4).Value = "https://s.technopoint.ru/thumb/st1/fit/800/650/afd6e1f91562be83338cd149f866811e/030a27e11b38b439b88a320f8122b836b9ec18b2de48232a1e0ef150ea491b32.jpg"

$ws.Cells.Item(8, 1).Value = "6.4`" Смартфон Huawei P40 Lite 128 ГБ черный"
$ws.Cells.Item(8, 2).Value = "'1640122"
$ws.Cells.Item(8, 3).Value = "'19999"
$ws.Cells.Item(8, 4).Value = "https://s.technopoint.ru/thumb/st4/fit/wm/800/650/3c4e73f9c97a92b220b0453802280290/59e88fa2e38f2de590507ca98480285c6f63528cd515a481861ed9050bf40451.jpg"

$ws.Cells.Item(9, 1).Value = "6.7`" Смартфон Samsung Galaxy S20+ 128 ГБ черный"
$ws.Cells.Item(9, 2).Value = "'1627804"
$ws.Cells.Item(9, 3).Value = "'79999"
$ws.Cells.Item(9, 4).Value = "https://s.technopoint.ru/thumb/st4/fit/800/650/53ccf05d955eb2daa680e5d5fe8e88db/1ebc2b0d87aa055e152648579180d03367ade9857a1176b8905b958971e13f3c.jpg"

$ws.Cells.Item(10, 1).Value = "5`" Смартфон bright & quick BQ 5016G CHOICE 16 ГБ красный"
$ws.Cells.Item(10, 2).Value = "'1646107"
$ws.Cells.Item(10, 3).Value = "'3999"
$ws.Cells.Item(10, 4).Value = "https://s.technopoint.ru/thumb/st1/fit/800/650/853093d89e3ff3610be956c1fd95c3e2/5c126190673c5232f844241ac215a576c667e82acb9852679508270d9ae99038.jpg"

$ws.Cells.Item(11, 1).Value = "5`" Смартфон bright & quick BQ 5016G CHOICE 16 ГБ черный"
$ws.Cells.Item(11, 2).Value = "'1646104"
$ws.Cells.Item(11, 3).Value = "'3999"
$ws.Cells.Item(11, 4).Value = "None"
